$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the header in C1: remove the trailing space from "Vendor/Merchant "
$ws.Range("C1").Value = "Vendor/Merchant"

# F5:G5 carried a redundant duplicate cell style; copying formats from the
# equivalent already-used style (E4) collapses it back to the shared one.
$ws.Range("E4").Copy()
$ws.Range("F5:G5").PasteSpecial(-4122)

# Move the active selection to C2 (as left after editing C1 and pressing Enter)
$ws.Range("C2").Select()
